$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HoReCa Bar Tavern_Night Club")

# Delete the obsolete KPI row (old row 37: "local 21" / ACTIVATION_OTHER / ANY OTHER ACTIVATION)
$ws.Rows(37).Delete()

# Fill in newly-populated brand cells on rows 28 and 30
$ws.Range("N28").Value = "Schweppes, Rich, Burn"
$ws.Range("N30").Value = "Schweppes, Coca-Cola, Burn"
